{"js": "// Office.js (Word JavaScript API) script\n// Applies the edits described by the commit:\n//  - Retitle the test plan from the \"Mon Budget\" app to \"cicd-todo-app\"\n//  - Update technology versions (DB, backend, frontend, app version)\n//  - Simplify the \"Environnement\" line for unit tests (drop \"/ Int\u00e9gration\")\n//  - Rewrite TC_A001 acceptance-test case and its matching UAT table row\n\nconst body = context.document.body;\n\n// 1) Title: 'Plan de Test \u2013 Application \"Mon Budget\"' -> 'Plan de Test \u2013 cicd-todo-app'\nlet results = body.search('Application \"Mon Budget\"', { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"cicd-todo-app\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 1b) Center the title paragraph\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.alignment = Word.Alignment.centered;\nawait context.sync();\n\n// 2) Objective paragraph: 'l'application \"Mon Budget \" \u00e0 travers' -> 'l'application \"todo-app\" \u00e0 travers'\nresults = body.search('\"Mon Budget \" \u00e0 travers', { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText('\"todo-app\" \u00e0 travers', Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Database version: MongoDB 8.2 -> MySQL\nresults = body.search(\"Base de donn\u00e9es : MongoDB 8.2\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Base de donn\u00e9es : MySQL\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Backend version: Node.js v18 -> Node.js v21.6.1\nresults = body.search(\"Backend : Node.js v18\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Backend : Node.js v21.6.1\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 5) Frontend framework: \"?\" -> \"vue 3.5.13\" (leading nbsp kept)\nresults = body.search(\"\\u00a0?\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"\\u00a0vue 3.5.13\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 6) App version: 1.0.0 -> 1.0.6\nresults = body.search(\"Version de l\\u2019application : 1.0.0\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Version de l\\u2019application : 1.0.6\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 7) Unit-test environment: \" DEV / Int\u00e9gration\" -> \" DEV\"\nresults = body.search(\" DEV / Int\u00e9gration\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\" DEV\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 8) TC_A001 acceptance test case text\nresults = body.search(\n  \"TC_A001 : L\\u2019utilisateur peut g\\u00e9rer ses d\\u00e9penses de mani\\u00e8re coh\\u00e9rente\",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\n    \"TC_A001 : Un retour \\u00e0 la ligne doit \\u00eatre fait si la description d\\u2019une t\\u00e2che est trop longue.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 9) UAT table row: feature description cell\nresults = body.search(\"Gestion des t\\u00e2ches sans formation\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Retour \\u00e0 la ligne, si description trop longue.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 10) UAT table row: expected result cell\nresults = body.search(\"L'utilisateur r\\u00e9ussit \\u00e0 g\\u00e9rer ses t\\u00e2ches\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"L'utilisateur voit les t\\u00e2ches sur plusieurs lignes\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the edits described by the commit:\n#  - Retitle the test plan from the \"Mon Budget\" app to \"cicd-todo-app\"\n#  - Update technology versions (DB, backend, frontend, app version)\n#  - Simplify the \"Environnement\" line for unit tests (drop \"/ Integration\")\n#  - Rewrite TC_A001 acceptance-test case and its matching UAT table row\n\n$d = $word.ActiveDocument\n\n# Use Find to locate the target text, then assign Range.Text directly\n# instead of Find.Execute(..., Replace:=wdReplaceAll, ...) so that straight\n# quotes/apostrophes in the replacement are kept verbatim (Word's Find\n# replacement path runs the text through AutoCorrect/\"smart quotes\").\nfunction Replace-Text($doc, $findText, $replaceText) {\n    $range = $doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $found = $find.Execute()\n    if ($found) {\n        $range.Text = $replaceText\n    }\n    return $found\n}\n\n# 1) Title: 'Plan de Test - Application \"Mon Budget\"' -> 'Plan de Test - cicd-todo-app'\nReplace-Text $d 'Application \"Mon Budget\"' 'cicd-todo-app'\n\n# 1b) Center the title paragraph (wdAlignParagraphCenter = 1)\n$d.Paragraphs(1).Alignment = 1\n\n# 2) Objective paragraph: l'application \"Mon Budget \" \u00e0 travers -> l'application \"todo-app\" \u00e0 travers\nReplace-Text $d '\"Mon Budget \" \u00e0 travers' '\"todo-app\" \u00e0 travers'\n\n# 3) Database version: MongoDB 8.2 -> MySQL\nReplace-Text $d 'Base de donn\u00e9es : MongoDB 8.2' 'Base de donn\u00e9es : MySQL'\n\n# 4) Backend version: Node.js v18 -> Node.js v21.6.1\nReplace-Text $d 'Backend : Node.js v18' 'Backend : Node.js v21.6.1'\n\n# 5) Frontend framework: \"?\" -> \"vue 3.5.13\" (leading nbsp kept)\n$nbsp = [char]0x00a0\n$frontendFind = $nbsp + '?'\n$frontendReplace = $nbsp + 'vue 3.5.13'\nReplace-Text $d $frontendFind $frontendReplace\n\n# 6) App version: 1.0.0 -> 1.0.6\nReplace-Text $d 'Version de l\u2019application : 1.0.0' 'Version de l\u2019application : 1.0.6'\n\n# 7) Unit-test environment: \" DEV / Int\u00e9gration\" -> \" DEV\"\nReplace-Text $d ' DEV / Int\u00e9gration' ' DEV'\n\n# 8) TC_A001 acceptance test case text\nReplace-Text $d 'TC_A001 : L\u2019utilisateur peut g\u00e9rer ses d\u00e9penses de mani\u00e8re coh\u00e9rente' 'TC_A001 : Un retour \u00e0 la ligne doit \u00eatre fait si la description d\u2019une t\u00e2che est trop longue.'\n\n# 9) UAT table row: feature description cell\nReplace-Text $d 'Gestion des t\u00e2ches sans formation' 'Retour \u00e0 la ligne, si description trop longue.'\n\n# 10) UAT table row: expected result cell\nReplace-Text $d \"L'utilisateur r\u00e9ussit \u00e0 g\u00e9rer ses t\u00e2ches\" \"L'utilisateur voit les t\u00e2ches sur plusieurs lignes\"\n"}
